# regen sval data to filter save games
#
# The per-player per-game stat columns (TB, d2S, K, IP) are recomputed from
# an updated save-game filter; each distinct old value for a column maps to
# exactly one new value. "sum" (G) is just TB + d2S + K + IP recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$epsilon = 0.0000001

function Map-Value($val, $pairs) {
    foreach ($pair in $pairs) {
        if ([Math]::Abs($val - $pair[0]) -lt $epsilon) {
            return $pair[1]
        }
    }
    throw "No mapping found for value $val"
}

$mapB = @(
    , @(3.641759230980763, 3.182878228561681)
    , @(1.592038685284193, 0.7287194209349384)
)
$mapC = @(
    , @(1.329362116779562, 1.65323645889881)
    , @(0.5496097884205744, 0.3375848360084654)
)
$mapD = @(
    , @(0.9945002603303881, 3.082599426703578)
    , @(0.4743886036040816, 0.7127328510149897)
    , @(0.2168782717193853, 0.1529057820181812)
)
$mapE = @(
    , @(0.3401128002026628, 0.4998867070740569)
)

for ($r = 2; $r -le 14; $r++) {
    $bOld = $ws.Cells.Item($r, 2).Value2
    $cOld = $ws.Cells.Item($r, 3).Value2
    $dOld = $ws.Cells.Item($r, 4).Value2
    $eOld = $ws.Cells.Item($r, 5).Value2

    $bNew = Map-Value $bOld $mapB
    $cNew = Map-Value $cOld $mapC
    $dNew = Map-Value $dOld $mapD
    $eNew = Map-Value $eOld $mapE

    $ws.Cells.Item($r, 2).Value = $bNew
    $ws.Cells.Item($r, 3).Value = $cNew
    $ws.Cells.Item($r, 4).Value = $dNew
    $ws.Cells.Item($r, 5).Value = $eNew

    $ws.Cells.Item($r, 7).Value = $bNew + $cNew + $dNew + $eNew
}
